# DC-Colos.xlsx data refresh
# 1) Rows 55-58 (ZRH, LYS, BOD, SKP) get reordered: the ZRH row moves down
#    to the bottom of that 4-row block, the other three shift up one row.
# 2) Row 264 (CTU / Chengdu, China) is removed entirely; every row below it
#    (265-331) shifts up by one, so the table shrinks from 331 to 330 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rotate the 4-row block A55:H58 -------------------------------
$colCount = 8

$row55 = @()
$row56 = @()
$row57 = @()
$row58 = @()
for ($col = 1; $col -le $colCount; $col++) {
    $row55 += ,$ws.Cells.Item(55, $col).Value()
    $row56 += ,$ws.Cells.Item(56, $col).Value()
    $row57 += ,$ws.Cells.Item(57, $col).Value()
    $row58 += ,$ws.Cells.Item(58, $col).Value()
}

for ($col = 1; $col -le $colCount; $col++) {
    $ws.Cells.Item(55, $col).Value = $row56[$col - 1]
    $ws.Cells.Item(56, $col).Value = $row57[$col - 1]
    $ws.Cells.Item(57, $col).Value = $row58[$col - 1]
    $ws.Cells.Item(58, $col).Value = $row55[$col - 1]
}

# --- Step 2: delete row 264 (CTU / Chengdu) and shift everything else up --
$ws.Rows(264).Delete()
